$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7f912436f5e0>),
                (''model'',
                 AdaBoostClassifier(estimator=RandomForestClassifier(max_depth=2,
                                                                     max_features=''log2'',
                                                                     min_samples_leaf=4,
                                                                     min_samples_split=6,
                                                                     n_estimators=50,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])'
$ws.Range("C2").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7f91240aa880>, ''scaler'': StandardScaler(), ''model__n_estimators'': 5, ''model__estimator__n_estimators'': 50, ''model__estimator__min_samples_split'': 6, ''model__estimator__min_samples_leaf'': 4, ''model__estimator__max_features'': ''log2'', ''model__estimator__max_depth'': 2, ''model__estimator__class_weight'': None}'
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7f91c788fb50>),
                (''model'',
                 AdaBoostClassifier(estimator=RandomForestClassifier(max_depth=5,
                                                                     min_samples_leaf=3,
                                                                     n_estimators=50,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])'
$ws.Range("C3").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7f912418a970>, ''scaler'': StandardScaler(), ''model__n_estimators'': 5, ''model__estimator__n_estimators'': 50, ''model__estimator__min_samples_split'': 2, ''model__estimator__min_samples_leaf'': 3, ''model__estimator__max_features'': ''sqrt'', ''model__estimator__max_depth'': 5, ''model__estimator__class_weight'': None}'
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7f91c76e6400>),
                (''model'',
                 AdaBoostClassifier(estimator=RandomForestClassifier(max_depth=6,
                                                                     min_samples_split=5,
                                                                     n_estimators=10,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])'
$ws.Range("C4").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7f91243c6dc0>, ''scaler'': RobustScaler(), ''model__n_estimators'': 5, ''model__estimator__n_estimators'': 10, ''model__estimator__min_samples_split'': 5, ''model__estimator__min_samples_leaf'': 1, ''model__estimator__max_features'': ''sqrt'', ''model__estimator__max_depth'': 6, ''model__estimator__class_weight'': None}'
